$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new cell density data rows (32, 40, 63, 81) ---
# Row 32 (sample POC-A6)
$ws.Range("B32").Value = 20231125
$ws.Range("C32").Value = "LG"
$ws.Range("D32").Value = 1
$ws.Range("E32").Value = 439
$ws.Range("F32").Value = 462
$ws.Range("G32").Value = 477
$ws.Range("H32").Value = 554
$ws.Range("I32").Value = 513
$ws.Range("J32").Value = 476

# Row 40 (sample POC-A14)
$ws.Range("B40").Value = 20231125
$ws.Range("C40").Value = "LG"
$ws.Range("D40").Value = 2
$ws.Range("E40").Value = 206
$ws.Range("F40").Value = 196
$ws.Range("G40").Value = 293
$ws.Range("H40").Value = 225
$ws.Range("I40").Value = 215
$ws.Range("J40").Value = 222

# Row 63 (sample POR-A7)
$ws.Range("B63").Value = 20231125
$ws.Range("C63").Value = "LG"
$ws.Range("D63").Value = 2
$ws.Range("E63").Value = 143
$ws.Range("F63").Value = 109
$ws.Range("G63").Value = 123
$ws.Range("H63").Value = 136
$ws.Range("I63").Value = 102
$ws.Range("J63").Value = 100

# Row 81 (sample POR-R19)
$ws.Range("B81").Value = 20231125
$ws.Range("C81").Value = "LG"
$ws.Range("D81").Value = 2
$ws.Range("E81").Value = 180
$ws.Range("F81").Value = 213
$ws.Range("G81").Value = 156
$ws.Range("H81").Value = 172
$ws.Range("I81").Value = 173
$ws.Range("J81").Value = 178

# Recalculate so the cv formulas in column L no longer show #DIV/0!
$excel.Calculate()

# --- Update the saved view state (pane/selection) ---
$ws.Application.ActiveWindow.ScrollRow = 48
$ws.Range("B66").Select()
